# Auto-generated edit script applying the Ixion_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the rows changed upstream.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 670.28
$ws.Range("I6").Value = 234.9375
$ws.Range("J6").Value = 1444.2222
$ws.Range("K6").Value = 704.8125
$ws.Range("L6").Value = 4332.6666
$ws.Range("M6").Value = -592.8125
$ws.Range("N6").Value = -4556.6666

$ws.Range("H9").Value = 83.2
$ws.Range("I9").Value = 83.2
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 83.2
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 85.8
$ws.Range("N9").ClearContents()

$ws.Range("H12").Value = 362.625
$ws.Range("I12").Value = 316.83334
$ws.Range("K12").Value = 316.83334
$ws.Range("M12").Value = -146.83334

$ws.Range("H21").Value = 49007.6
$ws.Range("I21").Value = 54509.5
$ws.Range("J21").Value = 27000
$ws.Range("K21").Value = 54509.5
$ws.Range("L21").Value = 27000
$ws.Range("M21").Value = -54041.5
$ws.Range("N21").Value = -27936

$ws.Range("H23").Value = 49007.6
$ws.Range("I23").Value = 54509.5
$ws.Range("J23").Value = 27000
$ws.Range("K23").Value = 54509.5
$ws.Range("L23").Value = 27000
$ws.Range("M23").Value = -54275.5
$ws.Range("N23").Value = -27468

$ws.Range("H29").Value = 1780
$ws.Range("I29").Value = 1666.6666
$ws.Range("J29").Value = 1950
$ws.Range("K29").Value = 4999.9998
$ws.Range("L29").Value = 5850
$ws.Range("M29").Value = -4718.9998
$ws.Range("N29").Value = -6412

$ws.Range("H38").Value = 1593.9286
$ws.Range("J38").Value = 2982.7144
$ws.Range("L38").Value = 8948.143199999999
$ws.Range("N38").Value = -9692.143199999999

$ws.Range("H58").Value = 615
$ws.Range("I58").Value = 615
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1845
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1695
$ws.Range("N58").ClearContents()

$ws.Range("H87").Value = 28933.334
$ws.Range("J87").Value = 28933.334
$ws.Range("L87").Value = 28933.334
$ws.Range("N87").Value = -31429.334

$ws.Range("H90").Value = 28933.334
$ws.Range("J90").Value = 28933.334
$ws.Range("L90").Value = 86800.00199999999
$ws.Range("N90").Value = -99280.00199999999

$ws.Range("H116").Value = 7248.421
$ws.Range("J116").Value = 2664.2856
$ws.Range("L116").Value = 2664.2856
$ws.Range("N116").Value = -9548.285599999999

$ws.Range("H132").Value = 1029.5588
$ws.Range("I132").Value = 1010.1613
$ws.Range("K132").Value = 3030.4839
$ws.Range("M132").Value = -500.4839000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4285.857
$ws.Range("J2").Value = 3925.25
$ws.Range("L2").Value = 3925.25
$ws.Range("N2").Value = -4151.25

$ws.Range("H32").Value = 4150.3066
$ws.Range("I32").Value = 3382.2654
$ws.Range("K32").Value = 3382.2654
$ws.Range("M32").Value = -3095.2654

$ws.Range("H45").Value = 10484.818
$ws.Range("I45").Value = 12481.889
$ws.Range("K45").Value = 12481.889
$ws.Range("M45").Value = -12104.889

$ws.Range("H61").Value = 14011
$ws.Range("I61").Value = 15512.375
$ws.Range("K61").Value = 15512.375
$ws.Range("M61").Value = -15300.375

$ws.Range("H116").Value = 4285.857
$ws.Range("J116").Value = 3925.25
$ws.Range("L116").Value = 3925.25
$ws.Range("N116").Value = -8513.25

$ws.Range("H122").Value = 2138028
$ws.Range("I122").Value = 2850215
$ws.Range("K122").Value = 8550645
$ws.Range("M122").Value = -8548195

$ws.Range("H132").Value = 2812.2163
$ws.Range("I132").Value = 1201.45
$ws.Range("J132").Value = 4707.2354
$ws.Range("K132").Value = 3604.35
$ws.Range("L132").Value = 14121.7062
$ws.Range("M132").Value = -1074.35
$ws.Range("N132").Value = -19181.7062

$ws.Range("H136").Value = 14011
$ws.Range("I136").Value = 15512.375
$ws.Range("K136").Value = 46537.125
$ws.Range("M136").Value = -43987.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4285.857
$ws.Range("J3").Value = 3925.25
$ws.Range("L3").Value = 3925.25
$ws.Range("N3").Value = -4153.25

$ws.Range("H94").Value = 1864.1578
$ws.Range("I94").Value = 1160.9
$ws.Range("J94").Value = 2645.5557
$ws.Range("K94").Value = 1160.9
$ws.Range("L94").Value = 2645.5557
$ws.Range("M94").Value = -709.9000000000001
$ws.Range("N94").Value = -3547.5557

$ws.Range("H134").Value = 5915.7095
$ws.Range("I134").Value = 10411.308
$ws.Range("J134").Value = 2668.889
$ws.Range("K134").Value = 31233.924
$ws.Range("L134").Value = 8006.667
$ws.Range("M134").Value = -28698.924
$ws.Range("N134").Value = -13076.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 70015
$ws.Range("J21").Value = 70015
$ws.Range("L21").Value = 70015
$ws.Range("N21").Value = -70485

$ws.Range("H22").Value = 530
$ws.Range("I22").Value = 459.0909
$ws.Range("J22").Value = 725
$ws.Range("K22").Value = 459.0909
$ws.Range("L22").Value = 725
$ws.Range("M22").Value = -109.0909
$ws.Range("N22").Value = -1425

$ws.Range("H58").Value = 1441.7587
$ws.Range("I58").Value = 1122.0555
$ws.Range("J58").Value = 1964.909
$ws.Range("K58").Value = 1122.0555
$ws.Range("L58").Value = 1964.909
$ws.Range("M58").Value = -919.0554999999999
$ws.Range("N58").Value = -2370.909

$ws.Range("H99").Value = 4115.1055
$ws.Range("I99").Value = 2895.923
$ws.Range("J99").Value = 6756.6665
$ws.Range("K99").Value = 2895.923
$ws.Range("L99").Value = 6756.6665
$ws.Range("M99").Value = -1397.923
$ws.Range("N99").Value = -9752.666499999999

$ws.Range("H105").Value = 1260
$ws.Range("I105").Value = 1300
$ws.Range("J105").Value = 1233.3334
$ws.Range("K105").Value = 1300
$ws.Range("L105").Value = 1233.3334
$ws.Range("M105").Value = 447
$ws.Range("N105").Value = -4727.3334

$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -22350

$ws.Range("H122").Value = 1999.25
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 1999.4
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 5998.200000000001
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -10898.2

$ws.Range("H126").Value = 4115.1055
$ws.Range("I126").Value = 2895.923
$ws.Range("J126").Value = 6756.6665
$ws.Range("K126").Value = 8687.769
$ws.Range("L126").Value = 20269.9995
$ws.Range("M126").Value = -6217.769
$ws.Range("N126").Value = -25209.9995

$ws.Range("H132").Value = 3357.739
$ws.Range("I132").Value = 3369.2942
$ws.Range("J132").Value = 3325
$ws.Range("K132").Value = 10107.8826
$ws.Range("L132").Value = 9975
$ws.Range("M132").Value = -7577.882599999999
$ws.Range("N132").Value = -15035

$ws.Range("H136").Value = 1441.7587
$ws.Range("I136").Value = 1122.0555
$ws.Range("J136").Value = 1964.909
$ws.Range("K136").Value = 3366.1665
$ws.Range("L136").Value = 5894.727000000001
$ws.Range("M136").Value = -816.1664999999998
$ws.Range("N136").Value = -10994.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2287
$ws.Range("I68").Value = 440.66666
$ws.Range("J68").Value = 4133.3335
$ws.Range("K68").Value = 1321.99998
$ws.Range("L68").Value = 12400.0005
$ws.Range("M68").Value = -510.9999800000001
$ws.Range("N68").Value = -14022.0005

$ws.Range("H71").Value = 2287
$ws.Range("I71").Value = 440.66666
$ws.Range("J71").Value = 4133.3335
$ws.Range("K71").Value = 3965.99994
$ws.Range("L71").Value = 37200.0015
$ws.Range("M71").Value = 90.0000600000003
$ws.Range("N71").Value = -45312.0015

$ws.Range("H107").Value = 286400.56
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 334084
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1002252
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -1006092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5521.8125
$ws.Range("I132").Value = 6452.3
$ws.Range("J132").Value = 3971
$ws.Range("K132").Value = 19356.9
$ws.Range("L132").Value = 11913
$ws.Range("M132").Value = -16826.9
$ws.Range("N132").Value = -16973

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2387.8235
$ws.Range("I7").Value = 2089.3
$ws.Range("J7").Value = 2814.2856
$ws.Range("K7").Value = 2089.3
$ws.Range("L7").Value = 2814.2856
$ws.Range("M7").Value = -1977.3
$ws.Range("N7").Value = -3038.2856

$ws.Range("H16").Value = 651.8461
$ws.Range("I16").Value = 557.2
$ws.Range("J16").Value = 967.3333
$ws.Range("K16").Value = 557.2
$ws.Range("L16").Value = 967.3333
$ws.Range("M16").Value = -387.2
$ws.Range("N16").Value = -1307.3333

$ws.Range("H126").Value = 2387.8235
$ws.Range("I126").Value = 2089.3
$ws.Range("J126").Value = 2814.2856
$ws.Range("K126").Value = 6267.900000000001
$ws.Range("L126").Value = 8442.856800000001
$ws.Range("M126").Value = -3797.900000000001
$ws.Range("N126").Value = -13382.8568

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2447.3635
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2984.2
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2984.2
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5730.2

$ws.Range("H100").Value = 1117.8
$ws.Range("I100").Value = 647.25
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1294.5
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -753.5
$ws.Range("N100").Value = -7082

$ws.Range("H122").Value = 4199.4
$ws.Range("I122").Value = 3666.1667
$ws.Range("J122").Value = 4999.25
$ws.Range("K122").Value = 10998.5001
$ws.Range("L122").Value = 14997.75
$ws.Range("M122").Value = -8548.500100000001
$ws.Range("N122").Value = -19897.75

$ws.Range("H132").Value = 2256.2222
$ws.Range("I132").Value = 1533.7778
$ws.Range("J132").Value = 2978.6667
$ws.Range("K132").Value = 4601.3334
$ws.Range("L132").Value = 8936.000100000001
$ws.Range("M132").Value = -2071.3334
$ws.Range("N132").Value = -13996.0001

